$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '268.95'
$ws.Range("D2").Style = "Normal"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '23.17'
$ws.Range("D3").Style = "Normal"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '6.337'
$ws.Range("D4").Style = "Normal"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.06177'
$ws.Range("D5").Style = "Normal"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '3.643'
$ws.Range("D6").Style = "Normal"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '6.681'
$ws.Range("D7").Style = "Normal"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.385'
$ws.Range("D8").Style = "Normal"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.8293'
$ws.Range("D9").Style = "Normal"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.01373'
$ws.Range("D10").Style = "Normal"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.1606'
$ws.Range("D11").Style = "Normal"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08222'
$ws.Range("D12").Style = "Normal"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.03192'
$ws.Range("D14").Style = "Normal"

$ws.Range("B15").Value = 'ProBitToken'

$ws.Range("C15").Value = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.1238'
$ws.Range("D15").Style = "Normal"

$ws.Range("E15").Value = '14ProBitTokenPROB'

$ws.Range("B16").Value = 'BitMartToken'

$ws.Range("C16").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.09319'
$ws.Range("D16").Style = "Normal"

$ws.Range("E16").Value = '15BitMartTokenBMX'

$ws.Range("B17").Value = 'MCDex'

$ws.Range("C17").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.856'
$ws.Range("D17").Style = "Normal"

$ws.Range("E17").Value = '16MCDexMCB'

$ws.Range("B18").Value = 'BitForexToken'

$ws.Range("C18").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.001675'
$ws.Range("D18").Style = "Normal"

$ws.Range("E18").Value = '17BitForexTokenBF'

$ws.Range("B19").Value = 'CoinExToken'

$ws.Range("C19").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.04733'
$ws.Range("D19").Style = "Normal"

$ws.Range("E19").Value = '18CoinExTokenCET'

$ws.Range("B20").Value = 'TigerCash'

$ws.Range("C20").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.006387'
$ws.Range("D20").Style = "Normal"

$ws.Range("E20").Value = '19TigerCashTCH'

$ws.Range("B21").Value = 'HotbitToken'

$ws.Range("C21").Value = 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.005656'
$ws.Range("D21").Style = "Normal"

$ws.Range("E21").Value = '20HotbitTokenHTB'

$ws.Range("B22").Value = 'BitKan'

$ws.Range("C22").Value = 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.001078'
$ws.Range("D22").Style = "Normal"

$ws.Range("E22").Value = '21BitKanKAN'

$ws.Range("B23").Value = 'NitroEx'

$ws.Range("C23").Value = 'https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx'

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.0001500'
$ws.Range("D23").Style = "Normal"

$ws.Range("E23").Value = '22NitroExNTX'

$ws.Range("B24").Value = 'LEO'

$ws.Range("C24").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '3.721'
$ws.Range("D24").Style = "Normal"

$ws.Range("E24").Value = '23LEOLEO'

$ws.Range("B25").Value = 'BTSEToken'

$ws.Range("C25").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.412'
$ws.Range("D25").Style = "Normal"

$ws.Range("E25").Value = '24BTSETokenBTSE'

$ws.Range("B26").Value = 'BitpandaEcosystemToken'

$ws.Range("C26").Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.3339'
$ws.Range("D26").Style = "Normal"

$ws.Range("E26").Value = '25BitpandaEcosystemTokenBEST'

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0002703'
$ws.Range("D27").Style = "Normal"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.04696'
$ws.Range("D40").Style = "Normal"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.006928'
$ws.Range("D41").Style = "Normal"

$ws.Range("B42").Value = 'BKEXToken'

$ws.Range("C42").Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1157'
$ws.Range("D42").Style = "Normal"

$ws.Range("E42").Value = '41BKEXTokenBKK'

$ws.Range("B43").Value = 'CEJI'

$ws.Range("C43").Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.003299'
$ws.Range("D43").Style = "Normal"

$ws.Range("E43").Value = '42CEJICEJI'

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.01153'
$ws.Range("D44").Style = "Normal"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00006259'
$ws.Range("D45").Style = "Normal"

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0009898'
$ws.Range("D46").Style = "Normal"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.9198'
$ws.Range("D48").Style = "Normal"

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.002284'
$ws.Range("D49").Style = "Normal"

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.00001400'
$ws.Range("D50").Style = "Normal"

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.01240'
$ws.Range("D51").Style = "Normal"
